$d = $word.ActiveDocument

# --- Locate the paragraph that ends the "Change log - 0.6" block ---
# ("2 Change commented line to properly run the SQLPSX.psm1.")
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "2 Change commented line to properly run the SQLPSX.psm1*") {
        $targetIndex = $i
        break
    }
}

# --- The _GoBack bookmark currently sits at the end of that paragraph; it
#     needs to move to the end of the new "1. Found bug ..." paragraph, so
#     detach it from here first. ---
$b = $d.Bookmarks.Item("_GoBack")
$b.Delete()

# --- Two empty (bold) paragraphs originally followed the target paragraph.
#     Only one blank line survives in the result, so remove the second one
#     and keep the first (untouched, so it keeps its "no run" shape). ---
$secondBlank = $d.Paragraphs.Item($targetIndex + 2)
$secondBlank.Range.Delete()

# --- Insert two new paragraphs after the remaining blank paragraph: one for
#     the "Change log - 0.7 ..." heading line and one for the "1. Found bug
#     ..." bullet. ---
$firstBlank = $d.Paragraphs.Item($targetIndex + 1)
$firstBlank.Range.InsertParagraphAfter()
$pChangeLog = $d.Paragraphs.Item($targetIndex + 2)

$pChangeLog.Range.InsertParagraphAfter()
$pFoundBug = $d.Paragraphs.Item($targetIndex + 3)

# --- Fill in the text (exclude the trailing paragraph-mark char from the
#     write so the paragraph's own formatting/paragraph-mark is kept). ---
$rChangeLog = $d.Range($pChangeLog.Range.Start, $pChangeLog.Range.End - 1)
$rChangeLog.Text = "Change log - 0.7 - 05/10/2010 17:20 - Max Trinidad"

$pFoundBug = $d.Paragraphs.Item($targetIndex + 3)
$rFoundBug = $d.Range($pFoundBug.Range.Start, $pFoundBug.Range.End - 1)
# Type the real text plus one throwaway trailing character. Adding a
# zero-length bookmark exactly at "end of the last run in the document" is
# unreliable, so we keep a placeholder character after the insertion point
# while the bookmark is created, then delete it afterwards.
$rFoundBug.Text = "1. Found bug missing path to Windows\system32 PowerShell modules.X"

$pFoundBug = $d.Paragraphs.Item($targetIndex + 3)
$bookmarkPos = $pFoundBug.Range.End - 2
$rZero = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $rZero)

# --- Remove the throwaway placeholder character; the bookmark stays put. ---
$pFoundBug = $d.Paragraphs.Item($targetIndex + 3)
$xPos = $pFoundBug.Range.End - 2
$d.Range($xPos, $xPos + 1).Delete()
